$wb = $excel.ActiveWorkbook

# Sheet "Detalle_Pasos" - column K (EnCQR-LSTM) corrections
$ws1 = $wb.Worksheets.Item("Detalle_Pasos")

$ws1.Range("K2").Value = 3.378118867953183
$ws1.Range("K3").Value = 3.130778146408153
$ws1.Range("K4").Value = 3.584090938229371
$ws1.Range("K5").Value = 3.63921249041532
$ws1.Range("K6").Value = 3.609880350062108
$ws1.Range("K7").Value = 3.705518894271908
$ws1.Range("K8").Value = 3.730152199424272
$ws1.Range("K9").Value = 3.666130060929392
$ws1.Range("K10").Value = 3.667144009078122
$ws1.Range("K11").Value = 3.738117578702732
$ws1.Range("K12").Value = 3.683428767109539
$ws1.Range("K13").Value = 3.378243987824212
$ws1.Range("K14").Value = 4.519516751660177
$ws1.Range("K15").Value = 13.92376511200364
$ws1.Range("K16").Value = 9.685203277971254
$ws1.Range("K17").Value = 3.942145086725423
$ws1.Range("K18").Value = 3.40172745782657
$ws1.Range("K19").Value = 3.565318386833658
$ws1.Range("K20").Value = 3.089110762771786
$ws1.Range("K21").Value = 2.695026445706866
$ws1.Range("K22").Value = 2.95751578239483
$ws1.Range("K23").Value = 3.170131192362806
$ws1.Range("K24").Value = 3.179276765730073
$ws1.Range("K25").Value = 3.333213506255694

# Sheet "Reliability_Data" - column C (Empirical) corrections
$ws2 = $wb.Worksheets.Item("Reliability_Data")

$ws2.Range("C809").Value = 0.04166666666666666
$ws2.Range("C810").Value = 0.04166666666666666
$ws2.Range("C811").Value = 0.04166666666666666
$ws2.Range("C812").Value = 0.04166666666666666
$ws2.Range("C813").Value = 0.04166666666666666
$ws2.Range("C814").Value = 0.04166666666666666
$ws2.Range("C815").Value = 0.04166666666666666
$ws2.Range("C830").Value = 0.08333333333333333
$ws2.Range("C831").Value = 0.08333333333333333
$ws2.Range("C832").Value = 0.08333333333333333
$ws2.Range("C833").Value = 0.125
$ws2.Range("C834").Value = 0.125
$ws2.Range("C835").Value = 0.1666666666666667
$ws2.Range("C836").Value = 0.2083333333333333
$ws2.Range("C837").Value = 0.2083333333333333
$ws2.Range("C838").Value = 0.2083333333333333
$ws2.Range("C839").Value = 0.25
$ws2.Range("C840").Value = 0.25
$ws2.Range("C841").Value = 0.25
$ws2.Range("C842").Value = 0.2916666666666667
$ws2.Range("C843").Value = 0.4166666666666667
$ws2.Range("C845").Value = 0.4583333333333333
$ws2.Range("C848").Value = 0.8333333333333334
$ws2.Range("C849").Value = 0.9166666666666666
$ws2.Range("C861").Value = 0.9583333333333334
$ws2.Range("C862").Value = 0.9583333333333334
$ws2.Range("C863").Value = 0.9583333333333334
$ws2.Range("C864").Value = 0.9583333333333334
$ws2.Range("C865").Value = 0.9583333333333334
